# Add new worksheet "ODI Batting Extra" after "ODI Batting"
$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "ODI Batting Extra"

# --- Header row (row 1) ---
$newSheet.Range("A1").Value = "MATCH_CODE"
$newSheet.Range("B1").Value = "BATTING_POSITION"
$newSheet.Range("C1").Value = "NUM_4"
$newSheet.Range("D1").Value = "NUM_6"
$newSheet.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$newSheet.Range("F1").Value = "MAN_OF_MATCH"

# Style header row: bold, thin border all around, centered horizontal / top vertical
$headerRange = $newSheet.Range("A1:F1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# --- Data cells that must stay TEXT (not auto-converted to numbers/percentages) ---
# Force text storage via NumberFormat "@" before assigning, so values like
# "4483", "4", "2", "15.48%" etc. are kept as literal strings.
# (Multi-area ranges only apply formatting to the first area, so set each
# contiguous block individually.)
$newSheet.Range("A2:A4").NumberFormat = "@"
$newSheet.Range("C2:D3").NumberFormat = "@"
$newSheet.Range("E2:E3").NumberFormat = "@"
$newSheet.Range("F2:F4").NumberFormat = "@"

$newSheet.Range("A2").Value = "4483"
$newSheet.Range("C2").Value = "4"
$newSheet.Range("D2").Value = "2"
$newSheet.Range("E2").Value = "15.48%"
$newSheet.Range("F2").Value = "NO"

$newSheet.Range("A3").Value = "4484"
$newSheet.Range("C3").Value = "3"
$newSheet.Range("D3").Value = "0"
$newSheet.Range("E3").Value = "8.56%"
$newSheet.Range("F3").Value = "NO"

$newSheet.Range("A4").Value = "4486"
$newSheet.Range("F4").Value = "NO"

# BATTING_POSITION column stores real numbers
$newSheet.Range("B2").Value = 1
$newSheet.Range("B3").Value = 1

# Row 4 (B4:E4) stays present but blank - keep the text format so the cells
# persist in the sheet as empty text cells instead of being dropped.
$newSheet.Range("B4:E4").NumberFormat = "@"
$newSheet.Range("B4:E4").Value = ""

# Restore the original active sheet (the workbook was on the first sheet
# before this edit; adding a sheet shouldn't change that).
$wb.Worksheets.Item(1).Activate()
